$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("latest")

$ws.Range("B2").Value = -0.2626095943911887
$ws.Range("C2").Value = 1.949383735764183
$ws.Range("D2").Value = 16.52272044676487
$ws.Range("E2").Value = 4.064814933888734
$ws.Range("F2").Value = 4.147487754911329
$ws.Range("G2").Value = 23
$ws.Range("B3").Value = -0.01304390213147815
$ws.Range("C3").Value = 1.772198340847518
$ws.Range("D3").Value = 10.9418263788771
$ws.Range("E3").Value = 3.307843161166669
$ws.Range("F3").Value = 3.385659091639844
$ws.Range("G3").Value = 22
$ws.Range("B4").Value = -0.660492300833325
$ws.Range("C4").Value = 1.197293372571959
$ws.Range("D4").Value = 5.006510981504795
$ws.Range("E4").Value = 2.237523403565825
$ws.Range("F4").Value = 2.190610405148969
$ws.Range("G4").Value = 21
$ws.Range("B5").Value = -0.05106328621311902
$ws.Range("C5").Value = 0.7433681293355469
$ws.Range("D5").Value = 1.840802665129192
$ws.Range("E5").Value = 1.356761830657537
$ws.Range("F5").Value = 1.391022042252352
$ws.Range("G5").Value = 20
$ws.Range("B6").Value = -0.06395863908153489
$ws.Range("C6").Value = 0.8000673122416484
$ws.Range("D6").Value = 1.420101517000478
$ws.Range("E6").Value = 1.191680123607203
$ws.Range("F6").Value = 1.222570274741316
$ws.Range("G6").Value = 19
$ws.Range("B7").Value = -0.03738556241501667
$ws.Range("C7").Value = 0.71332745744033
$ws.Range("D7").Value = 0.9386152008072798
$ws.Range("E7").Value = 0.9688215526129049
$ws.Range("F7").Value = 0.9961666341101381
$ws.Range("G7").Value = 18
$ws.Range("B8").Value = 0.06091063799584202
$ws.Range("C8").Value = 0.6343285541743425
$ws.Range("D8").Value = 0.7182370180715076
$ws.Range("E8").Value = 0.8474886536535505
$ws.Range("F8").Value = 0.8713121394001694
$ws.Range("G8").Value = 17
$ws.Range("B9").Value = 0.1648078376396195
$ws.Range("C9").Value = 0.5511977509794913
$ws.Range("D9").Value = 0.4449526034813988
$ws.Range("E9").Value = 0.6670476770676882
$ws.Range("F9").Value = 0.6675655114989703
$ws.Range("G9").Value = 16
$ws.Range("B10").Value = 0.1934705323399399
$ws.Range("C10").Value = 0.5837592819091525
$ws.Range("D10").Value = 0.5633658838424578
$ws.Range("E10").Value = 0.7505770339162116
$ws.Range("F10").Value = 0.7506675864274018
$ws.Range("G10").Value = 15
$ws.Range("B11").Value = 0.2495182209949923
$ws.Range("C11").Value = 0.5857948382811914
$ws.Range("D11").Value = 0.4537182255294631
$ws.Range("E11").Value = 0.6735860936283224
$ws.Range("F11").Value = 0.6492850719707854
$ws.Range("G11").Value = 14
